$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.669886159284943
$ws.Cells.Item(2, 3).Value = 0.339145502435855
$ws.Cells.Item(2, 5).Value = 0.06909409277560119
$ws.Cells.Item(2, 6).Value = 2.969707982845279
$ws.Cells.Item(2, 7).Value = 0.002550694924269734
$ws.Cells.Item(2, 9).Value = 1.702634974721065
$ws.Cells.Item(2, 10).Value = 0.1171450160749998
$ws.Cells.Item(2, 12).Value = 0.4373352176423992
$ws.Cells.Item(2, 14).Value = 1.875600795357766

$ws.Cells.Item(3, 2).Value = 2.545457965338414
$ws.Cells.Item(3, 3).Value = 0.3026972707380367
$ws.Cells.Item(3, 5).Value = 0.0687517779157738
$ws.Cells.Item(3, 6).Value = 2.950025442101591
$ws.Cells.Item(3, 7).Value = 0.002556428290683568
$ws.Cells.Item(3, 9).Value = 1.699636447176772
$ws.Cells.Item(3, 10).Value = 0.1184482790675361
$ws.Cells.Item(3, 12).Value = 0.4307726991647485
$ws.Cells.Item(3, 14).Value = 1.897502744192515

$ws.Cells.Item(4, 2).Value = 2.470569751514063
$ws.Cells.Item(4, 3).Value = 0.2804674503497893
$ws.Cells.Item(4, 5).Value = 0.06853787779973564
$ws.Cells.Item(4, 6).Value = 2.939592481797831
$ws.Cells.Item(4, 7).Value = 0.002560133084547107
$ws.Cells.Item(4, 9).Value = 1.698771679272426
$ws.Cells.Item(4, 10).Value = 0.1193003671316113
$ws.Cells.Item(4, 12).Value = 0.4269370500407206
$ws.Cells.Item(4, 14).Value = 1.911642905357031

$ws.Cells.Item(5, 2).Value = 2.440430928523369
$ws.Cells.Item(5, 3).Value = 0.271445244896654
$ws.Cells.Item(5, 5).Value = 0.06844975248008667
$ws.Cells.Item(5, 6).Value = 2.935755375093819
$ws.Cells.Item(5, 7).Value = 0.002561689369569123
$ws.Cells.Item(5, 9).Value = 1.698663941636376
$ws.Cells.Item(5, 10).Value = 0.1196606056208402
$ws.Cells.Item(5, 12).Value = 0.4254227990302297
$ws.Cells.Item(5, 14).Value = 1.917578858329783

$ws.Cells.Item(6, 2).Value = 2.435449238199681
$ws.Cells.Item(6, 3).Value = 0.2699492983996095
$ws.Cells.Item(6, 5).Value = 0.06843506068993532
$ws.Cells.Item(6, 6).Value = 2.935143225421413
$ws.Cells.Item(6, 7).Value = 0.002561950605810112
$ws.Cells.Item(6, 9).Value = 1.698660804325911
$ws.Cells.Item(6, 10).Value = 0.1197212073069736
$ws.Cells.Item(6, 12).Value = 0.4251743091594022
$ws.Cells.Item(6, 14).Value = 1.918575000891707

$ws.Cells.Item(7, 2).Value = 2.470161757007816
$ws.Cells.Item(7, 3).Value = 0.2803456267067475
$ws.Cells.Item(7, 5).Value = 0.06853669322588085
$ws.Cells.Item(7, 6).Value = 2.939539056712775
$ws.Cells.Item(7, 7).Value = 0.002560153884400435
$ws.Cells.Item(7, 9).Value = 1.698769236792657
$ws.Cells.Item(7, 10).Value = 0.1193051728280174
$ws.Cells.Item(7, 12).Value = 0.4269164305884772
$ws.Cells.Item(7, 14).Value = 1.911722256790153

$ws.Cells.Item(8, 2).Value = 2.626668846991436
$ws.Cells.Item(8, 3).Value = 0.3265465135016257
$ws.Cells.Item(8, 5).Value = 0.06897681924079446
$ws.Cells.Item(8, 6).Value = 2.962577878627471
$ws.Cells.Item(8, 7).Value = 0.002552633597004084
$ws.Cells.Item(8, 9).Value = 1.701397925572792
$ws.Cells.Item(8, 10).Value = 0.1175835895918826
$ws.Cells.Item(8, 12).Value = 0.4350322919715239
$ws.Cells.Item(8, 14).Value = 1.883008692522825

$ws.Cells.Item(9, 2).Value = 2.94564744316051
$ws.Cells.Item(9, 3).Value = 0.4183817117385615
$ws.Cells.Item(9, 5).Value = 0.06981155157622609
$ws.Cells.Item(9, 6).Value = 3.020920508378495
$ws.Cells.Item(9, 7).Value = 0.002539342745531204
$ws.Cells.Item(9, 9).Value = 1.714340542354776
$ws.Cells.Item(9, 10).Value = 0.1146209836763568
$ws.Cells.Item(9, 12).Value = 0.4524826425562338
$ws.Cells.Item(9, 14).Value = 1.832209537410272

$ws.Cells.Item(10, 2).Value = 3.187496178796209
$ws.Cells.Item(10, 3).Value = 0.4866807572074094
$ws.Cells.Item(10, 5).Value = 0.07040917988853623
$ws.Cells.Item(10, 6).Value = 3.071894426536431
$ws.Cells.Item(10, 7).Value = 0.002530455465181472
$ws.Cells.Item(10, 9).Value = 1.728657823457695
$ws.Cells.Item(10, 10).Value = 0.1126988145636325
$ws.Cells.Item(10, 12).Value = 0.4662383686759171
$ws.Cells.Item(10, 14).Value = 1.798264960438942

$ws.Cells.Item(11, 2).Value = 3.299179155205707
$ws.Cells.Item(11, 3).Value = 0.5179481669991333
$ws.Cells.Item(11, 5).Value = 0.07067802078017227
$ws.Cells.Item(11, 6).Value = 3.096863914660076
$ws.Cells.Item(11, 7).Value = 0.002526600740376232
$ws.Cells.Item(11, 9).Value = 1.736228758892025
$ws.Cells.Item(11, 10).Value = 0.1118801617174405
$ws.Cells.Item(11, 12).Value = 0.4726992095948219
$ws.Cells.Item(11, 14).Value = 1.783560425042509

$ws.Cells.Item(12, 2).Value = 3.341712076502688
$ws.Cells.Item(12, 3).Value = 0.529817998637327
$ws.Cells.Item(12, 5).Value = 0.07077941847285896
$ws.Cells.Item(12, 6).Value = 3.106576765719268
$ws.Cells.Item(12, 7).Value = 0.002525167941860531
$ws.Cells.Item(12, 9).Value = 1.739248863894943
$ws.Cells.Item(12, 10).Value = 0.1115782241431091
$ws.Cells.Item(12, 12).Value = 0.4751749572812258
$ws.Cells.Item(12, 14).Value = 1.778098599643577

$ws.Cells.Item(13, 2).Value = 3.332541106062081
$ws.Cells.Item(13, 3).Value = 0.5272602865285876
$ws.Cells.Item(13, 5).Value = 0.07075759829080708
$ws.Cells.Item(13, 6).Value = 3.104473458203557
$ws.Cells.Item(13, 7).Value = 0.002525475326608974
$ws.Cells.Item(13, 9).Value = 1.738591602066634
$ws.Cells.Item(13, 10).Value = 0.1116428921218073
$ws.Cells.Item(13, 12).Value = 0.4746404645033806
$ws.Cells.Item(13, 14).Value = 1.779270158747885

$ws.Cells.Item(14, 2).Value = 3.302673524792681
$ws.Cells.Item(14, 3).Value = 0.5189241067977832
$ws.Cells.Item(14, 5).Value = 0.07068637081305562
$ws.Cells.Item(14, 6).Value = 3.09765782951186
$ws.Cells.Item(14, 7).Value = 0.002526482325062175
$ws.Cells.Item(14, 9).Value = 1.736474149337127
$ws.Cells.Item(14, 10).Value = 0.1118551591279733
$ws.Cells.Item(14, 12).Value = 0.4729023066528981
$ws.Cells.Item(14, 14).Value = 1.783108940921146

$ws.Cells.Item(15, 2).Value = 3.284410203099355
$ws.Cells.Item(15, 3).Value = 0.5138218360210658
$ws.Cells.Item(15, 5).Value = 0.07064268983254474
$ws.Cells.Item(15, 6).Value = 3.093516626860151
$ws.Cells.Item(15, 7).Value = 0.002527102639581165
$ws.Cells.Item(15, 9).Value = 1.735197124576175
$ws.Cells.Item(15, 10).Value = 0.1119862311006266
$ws.Cells.Item(15, 12).Value = 0.4718414313479968
$ws.Cells.Item(15, 14).Value = 1.785474184873856

$ws.Cells.Item(16, 2).Value = 3.18023106952171
$ws.Cells.Item(16, 3).Value = 0.4846414324914576
$ws.Cells.Item(16, 5).Value = 0.07039155206869596
$ws.Cells.Item(16, 6).Value = 3.070298565749255
$ws.Cells.Item(16, 7).Value = 0.002530711153316976
$ws.Cells.Item(16, 9).Value = 1.728184420731779
$ws.Cells.Item(16, 10).Value = 0.1127534416926075
$ws.Cells.Item(16, 12).Value = 0.4658202225956245
$ws.Cells.Item(16, 14).Value = 1.79924080801716

$ws.Cells.Item(17, 2).Value = 3.116748062477086
$ws.Cells.Item(17, 3).Value = 0.4667916430051378
$ws.Cells.Item(17, 5).Value = 0.07023673277131248
$ws.Cells.Item(17, 6).Value = 3.056512171577509
$ws.Cells.Item(17, 7).Value = 0.002532972939032973
$ws.Cells.Item(17, 9).Value = 1.724154041476027
$ws.Cells.Item(17, 10).Value = 0.1132384161357542
$ws.Cells.Item(17, 12).Value = 0.4621784262900093
$ws.Cells.Item(17, 14).Value = 1.80787529327614

$ws.Cells.Item(18, 2).Value = 3.080390912804262
$ws.Cells.Item(18, 3).Value = 0.456543470639815
$ws.Cells.Item(18, 5).Value = 0.07014739852522123
$ws.Cells.Item(18, 6).Value = 3.048750176967189
$ws.Cells.Item(18, 7).Value = 0.002534291576238147
$ws.Cells.Item(18, 9).Value = 1.721935389229486
$ws.Cells.Item(18, 10).Value = 0.1135226050625491
$ws.Cells.Item(18, 12).Value = 0.4601028999283159
$ws.Cells.Item(18, 14).Value = 1.812910963665445

$ws.Cells.Item(19, 2).Value = 3.068107854099935
$ws.Cells.Item(19, 3).Value = 0.4530767692308473
$ws.Cells.Item(19, 5).Value = 0.07011710146490868
$ws.Cells.Item(19, 6).Value = 3.046150840422143
$ws.Cells.Item(19, 7).Value = 0.002534741092222197
$ws.Cells.Item(19, 9).Value = 1.721201252505665
$ws.Cells.Item(19, 10).Value = 0.1136197257379337
$ws.Cells.Item(19, 12).Value = 0.4594034514823591
$ws.Cells.Item(19, 14).Value = 1.81462785034207

$ws.Cells.Item(20, 2).Value = 3.123489718276346
$ws.Cells.Item(20, 3).Value = 0.4686898553793526
$ws.Cells.Item(20, 5).Value = 0.07025324297621172
$ws.Cells.Item(20, 6).Value = 3.057962403995987
$ws.Cells.Item(20, 7).Value = 0.002532730335483624
$ws.Cells.Item(20, 9).Value = 1.72457277538868
$ws.Cells.Item(20, 10).Value = 0.1131862466282563
$ws.Cells.Item(20, 12).Value = 0.4625641214429521
$ws.Cells.Item(20, 14).Value = 1.806948957414427

$ws.Cells.Item(21, 2).Value = 3.311439803597068
$ws.Cells.Item(21, 3).Value = 0.5213718342920401
$ws.Cells.Item(21, 5).Value = 0.07070730285500781
$ws.Cells.Item(21, 6).Value = 3.099652747216226
$ws.Cells.Item(21, 7).Value = 0.002526185816370847
$ws.Cells.Item(21, 9).Value = 1.737091931658441
$ws.Cells.Item(21, 10).Value = 0.1117925917706124
$ws.Cells.Item(21, 12).Value = 0.4734120546589367
$ws.Cells.Item(21, 14).Value = 1.781978504294734

$ws.Cells.Item(22, 2).Value = 3.435681578158778
$ws.Cells.Item(22, 3).Value = 0.5559752665052429
$ws.Cells.Item(22, 5).Value = 0.07100170122226457
$ws.Cells.Item(22, 6).Value = 3.128401026995789
$ws.Cells.Item(22, 7).Value = 0.00252206532973383
$ws.Cells.Item(22, 9).Value = 1.746167122466531
$ws.Cells.Item(22, 10).Value = 0.1109288117704352
$ws.Cells.Item(22, 12).Value = 0.4806717880974389
$ws.Cells.Item(22, 14).Value = 1.766279646994704

$ws.Cells.Item(23, 2).Value = 3.36924225322042
$ws.Cells.Item(23, 3).Value = 0.5374905944614738
$ws.Cells.Item(23, 5).Value = 0.07084478138264583
$ws.Cells.Item(23, 6).Value = 3.112919725704387
$ws.Cells.Item(23, 7).Value = 0.002524250220088191
$ws.Cells.Item(23, 9).Value = 1.74124145126018
$ws.Cells.Item(23, 10).Value = 0.1113855046695704
$ws.Cells.Item(23, 12).Value = 0.4767815992557303
$ws.Cells.Item(23, 14).Value = 1.77460146039693

$ws.Cells.Item(24, 2).Value = 3.12044137950204
$ws.Cells.Item(24, 3).Value = 0.467831630323758
$ws.Cells.Item(24, 5).Value = 0.07024577973187229
$ws.Cells.Item(24, 6).Value = 3.057306243235445
$ws.Cells.Item(24, 7).Value = 0.002532839959689726
$ws.Cells.Item(24, 9).Value = 1.724383159163182
$ws.Cells.Item(24, 10).Value = 0.1132098157378962
$ws.Cells.Item(24, 12).Value = 0.4623896919991495
$ws.Cells.Item(24, 14).Value = 1.807367530907484

$ws.Cells.Item(25, 2).Value = 2.858049187585209
$ws.Cells.Item(25, 3).Value = 0.3933979423747473
$ws.Cells.Item(25, 5).Value = 0.06958860258011201
$ws.Cells.Item(25, 6).Value = 3.003719357513447
$ws.Cells.Item(25, 7).Value = 0.002542783427391334
$ws.Cells.Item(25, 9).Value = 1.709999265556718
$ws.Cells.Item(25, 10).Value = 0.1153779409122997
$ws.Cells.Item(25, 12).Value = 0.4475976168847637
$ws.Cells.Item(25, 14).Value = 1.845360375772287
